# Regenerate the localization-status report for a new source file
# (919668ce-... -> 9f87e1ae-...), refreshed handoff hashes/timestamps,
# and a cleared "Latest Target File" / "Latest Handback File" for both
# locales, as described in the "Generate Report for Handoff" commit.

$wb = $excel.ActiveWorkbook

$newGuid = "9f87e1ae-a573-43d9-88b3-226cba314a9e"
$newHash = "e35c12d45cf2407da70cd7457e7b874b8675f7ad"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------
# Capture existing hyperlink addresses (must iterate with foreach;
# indexed .Item() access on this engine returns a stale/empty object).
# ---------------------------------------------------------------
$ws1LinkAddr = $null
foreach ($h in $ws1.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$2') { $ws1LinkAddr = $h.Address() }
}

$ws2LinkAddr = $null
foreach ($h in $ws2.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') { $ws2LinkAddr = $h.Address() }
}

$ws3LinkAddr = $null
foreach ($h in $ws3.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') { $ws3LinkAddr = $h.Address() }
}

# ---------------------------------------------------------------
# Drop every hyperlink on each sheet (this engine's Hyperlinks.Delete()
# clears the whole sheet regardless of which Range it was called from),
# then re-create only the ones that should survive, with fresh text.
# ---------------------------------------------------------------
$ws1.Range("A1").Hyperlinks.Delete()
$ws2.Range("A1").Hyperlinks.Delete()
$ws3.Range("A1").Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("B2"), $ws1LinkAddr, [Type]::Missing, [Type]::Missing, ("e2e\" + $newGuid + ".md"))
$ws2.Hyperlinks.Add($ws2.Range("A2"), $ws2LinkAddr, [Type]::Missing, [Type]::Missing, ($newGuid + ".md"))
$ws3.Hyperlinks.Add($ws3.Range("A2"), $ws3LinkAddr, [Type]::Missing, [Type]::Missing, ($newGuid + ".md"))

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws1.Range("A2").Value = $newGuid + ".md"
$ws1.Range("B2").Value = "e2e\" + $newGuid + ".md"
$ws1.Range("G2").Value = "2016-08-27 19:06:32"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws2.Range("A2").Value = $newGuid + ".md"
$ws2.Range("G2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-27 19:06:28"
$ws2.Range("I2").Value = ""
$ws2.Range("I2").Style = "Normal"
$ws2.Range("J2").Value = ""
$ws2.Range("J2").Style = "Normal"
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Columns.Item(9).ColumnWidth = 18.6506053379604 - 5/6
$ws2.Columns.Item(10).ColumnWidth = 21.7054770333426 - 5/6

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws3.Range("A2").Value = $newGuid + ".md"
$ws3.Range("G2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$ws3.Range("H2").Value = "2016-08-27 19:06:32"
$ws3.Range("I2").Value = ""
$ws3.Range("I2").Style = "Normal"
$ws3.Range("J2").Value = ""
$ws3.Range("J2").Style = "Normal"
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Columns.Item(9).ColumnWidth = 18.6506053379604 - 5/6
$ws3.Columns.Item(10).ColumnWidth = 21.7054770333426 - 5/6
